$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L43").Value = 3840.1667
$ws.Range("J43").Value = 3840.1667
$ws.Range("K43").Value = 2500
$ws.Range("H43").Value = 3648.7144
$ws.Range("I43").Value = 2500
$ws.Range("M43").Value = -2431
$ws.Range("N43").Value = -3978.1667
$ws.Range("L74").Value = 41666.668
$ws.Range("J74").Value = 41666.668
$ws.Range("K74").Value = 6674.9165
$ws.Range("H74").Value = 13673.267
$ws.Range("I74").Value = 6674.9165
$ws.Range("M74").Value = -5738.9165
$ws.Range("N74").Value = -43538.668
$ws.Range("L77").Value = 208333.34
$ws.Range("J77").Value = 41666.668
$ws.Range("K77").Value = 33374.5825
$ws.Range("H77").Value = 13673.267
$ws.Range("I77").Value = 6674.9165
$ws.Range("M77").Value = -28694.5825
$ws.Range("N77").Value = -217693.34
$ws.Range("L98").Value = 5000
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 3804
$ws.Range("H98").Value = 4103
$ws.Range("I98").Value = 3804
$ws.Range("M98").Value = -2306
$ws.Range("N98").Value = -7996
$ws.Range("L122").Value = 15000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11412
$ws.Range("H122").Value = 4103
$ws.Range("I122").Value = 3804
$ws.Range("M122").Value = -8962
$ws.Range("N122").Value = -19900
$ws.Range("K137").Value = 131743.254
$ws.Range("H137").Value = 38319.43
$ws.Range("I137").Value = 43914.418
$ws.Range("M137").Value = -129193.254

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 21138.92
$ws.Range("H32").Value = 19324.018
$ws.Range("I32").Value = 21138.92
$ws.Range("M32").Value = -20851.92
$ws.Range("L45").Value = 3681.5217
$ws.Range("J45").Value = 3681.5217
$ws.Range("K45").Value = 1613.8334
$ws.Range("H45").Value = 3253.724
$ws.Range("I45").Value = 1613.8334
$ws.Range("M45").Value = -1236.8334
$ws.Range("N45").Value = -4435.521699999999
$ws.Range("K61").Value = 3839.7727
$ws.Range("H61").Value = 3798.8667
$ws.Range("I61").Value = 3839.7727
$ws.Range("M61").Value = -3627.7727
$ws.Range("L122").Value = 3000
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4326.428400000001
$ws.Range("H122").Value = 1426.8966
$ws.Range("I122").Value = 1442.1428
$ws.Range("M122").Value = -1876.428400000001
$ws.Range("N122").Value = -7900
$ws.Range("L132").Value = 8097.999899999999
$ws.Range("J132").Value = 2699.3333
$ws.Range("K132").Value = 101530.032
$ws.Range("H132").Value = 28925.87
$ws.Range("I132").Value = 33843.344
$ws.Range("M132").Value = -99000.03199999999
$ws.Range("N132").Value = -13157.9999
$ws.Range("K136").Value = 11519.3181
$ws.Range("H136").Value = 3798.8667
$ws.Range("I136").Value = 3839.7727
$ws.Range("M136").Value = -8969.3181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K107").Value = 3000
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("M107").Value = -1080
$ws.Range("L134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8049.299999999999
$ws.Range("H134").Value = 2683.1
$ws.Range("I134").Value = 2683.1
$ws.Range("M134").Value = -5514.299999999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L22").Value = 2324.875
$ws.Range("J22").Value = 2324.875
$ws.Range("K22").Value = 1095.6
$ws.Range("H22").Value = 1852.0769
$ws.Range("I22").Value = 1095.6
$ws.Range("M22").Value = -745.5999999999999
$ws.Range("N22").Value = -3024.875
$ws.Range("L31").Value = 4799.4
$ws.Range("J31").Value = 4799.4
$ws.Range("K31").Value = 1970.8572
$ws.Range("H31").Value = 3149.4167
$ws.Range("I31").Value = 1970.8572
$ws.Range("M31").Value = -1675.8572
$ws.Range("N31").Value = -5389.4
$ws.Range("L34").Value = 4799.4
$ws.Range("J34").Value = 4799.4
$ws.Range("K34").Value = 1970.8572
$ws.Range("H34").Value = 3149.4167
$ws.Range("I34").Value = 1970.8572
$ws.Range("M34").Value = -1768.8572
$ws.Range("N34").Value = -5203.4
$ws.Range("L132").Value = 8741.849999999999
$ws.Range("J132").Value = 2913.95
$ws.Range("K132").Value = 4903.9092
$ws.Range("H132").Value = 2460
$ws.Range("I132").Value = 1634.6364
$ws.Range("M132").Value = -2373.9092
$ws.Range("N132").Value = -13801.85
$ws.Range("K134").Value = 339441
$ws.Range("H134").Value = 93147
$ws.Range("I134").Value = 113147
$ws.Range("M134").Value = -336906

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K2").Value = 54.416668
$ws.Range("H2").Value = 100.388885
$ws.Range("I2").Value = 54.416668
$ws.Range("M2").Value = 58.583332
$ws.Range("K102").Value = 5095.4443
$ws.Range("H102").Value = 5095.4443
$ws.Range("I102").Value = 5095.4443
$ws.Range("M102").Value = -3473.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 2410.3
$ws.Range("H7").Value = 2554.818
$ws.Range("I7").Value = 2410.3
$ws.Range("M7").Value = -2298.3
$ws.Range("L16").Value = 8000
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 8394.799999999999
$ws.Range("H16").Value = 8329
$ws.Range("I16").Value = 8394.799999999999
$ws.Range("M16").Value = -8224.799999999999
$ws.Range("N16").Value = -8340
$ws.Range("L46").Value = 4793.4443
$ws.Range("J46").Value = 4793.4443
$ws.Range("K46").Value = 27621.777
$ws.Range("H46").Value = 16207.611
$ws.Range("I46").Value = 27621.777
$ws.Range("M46").Value = -27433.777
$ws.Range("N46").Value = -5169.4443
$ws.Range("L55").Value = 1212.5
$ws.Range("J55").Value = 1212.5
$ws.Range("K55").Value = 142.23077
$ws.Range("H55").Value = 394.05884
$ws.Range("I55").Value = 142.23077
$ws.Range("M55").Value = 30.76922999999999
$ws.Range("N55").Value = -1558.5
$ws.Range("K61").Value = 1398.1
$ws.Range("H61").Value = 2141.8572
$ws.Range("I61").Value = 1398.1
$ws.Range("M61").Value = -1196.1
$ws.Range("L68").Value = 4800
$ws.Range("J68").Value = 4800
$ws.Range("K68").Value = 5873.25
$ws.Range("H68").Value = 5658.6
$ws.Range("I68").Value = 5873.25
$ws.Range("M68").Value = -5124.25
$ws.Range("N68").Value = -6298
$ws.Range("L71").Value = 24000
$ws.Range("J71").Value = 4800
$ws.Range("K71").Value = 29366.25
$ws.Range("H71").Value = 5658.6
$ws.Range("I71").Value = 5873.25
$ws.Range("M71").Value = -25622.25
$ws.Range("N71").Value = -31488
$ws.Range("K113").Value = 1398.1
$ws.Range("H113").Value = 2141.8572
$ws.Range("I113").Value = 1398.1
$ws.Range("M113").Value = 771.9000000000001
$ws.Range("K126").Value = 7230.900000000001
$ws.Range("H126").Value = 2554.818
$ws.Range("I126").Value = 2410.3
$ws.Range("M126").Value = -4760.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L100").Value = 1978
$ws.Range("J100").Value = 989
$ws.Range("K100").Value = 1980
$ws.Range("H100").Value = 989.5
$ws.Range("I100").Value = 990
$ws.Range("M100").Value = -1439
$ws.Range("N100").Value = -3060
$ws.Range("L126").Value = 26556
$ws.Range("J126").Value = 8852
$ws.Range("H126").Value = 174122.33
$ws.Range("N126").Value = -31496
